# Implemented GEP with constant region mutation and inversion
# Adds two new result columns (G:J) for two additional GEP parameter runs:
#   GEP(pop25, gen40) -> columns G (AVG) / H (Gen)
#   GEP(pop40, gen25) -> columns I (AVG) / J (Gen)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the bordered header style (same style used by
# the existing B24/C24/E24 "n_perfect solution" row) onto the new header
# cells, then set their labels.
[void]$ws.Range("B24").Copy()
$ws.Range("G1:J1").PasteSpecial(-4122)
$ws.Range("G1").Value = "GEP(pop25, gen40)"
$ws.Range("H1").Value = "Gen"
$ws.Range("I1").Value = "GEP(pop40, gen25)"
$ws.Range("J1").Value = "Gen"

# --- Data rows: AVG (G / I) and Gen (H / J) values for the two new runs.
$ws.Range("G2").Value = 997.76899699820694
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 10

$ws.Range("G3").Value = 999.27344704546101
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 12

$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 998.14113464749403

$ws.Range("G5").Value = 997.76899699820694
$ws.Range("I5").Value = 997.76899699820694

$ws.Range("I6").Value = 999.79219766162396

$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 7

$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 9

$ws.Range("I9").Value = 997.54305803694604

$ws.Range("I10").Value = 998.30575172668398

$ws.Range("I11").Value = 997.76899699820694

$ws.Range("I12").Value = 999.92484796702797

# --- Column widths for the two new value columns (match column E's width).
$ws.Columns("G").ColumnWidth = 17.022135416666668
$ws.Columns("I").ColumnWidth = 17.022135416666668

# --- Conditional formatting: extend the color-scale coverage to the new
# header cells (G1:J1), matching the existing rule used on B1:C24/E1:E24.
$csNew = $ws.Range("G1:J1").FormatConditions.AddColorScale(3)

# --- Selection moved by the editor while working in column C.
[void]$ws.Range("C10").Select()
